$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 203.45454
$ws.Range("I55").Value = 70.888885
$ws.Range("J55").Value = 800
$ws.Range("K55").Value = 70.888885
$ws.Range("L55").Value = 800
$ws.Range("M55").Value = 143.111115
$ws.Range("N55").Value = -1228
$ws.Range("H116").Value = 2892.3076
$ws.Range("I116").Value = 2660
$ws.Range("J116").Value = 3666.6667
$ws.Range("K116").Value = 2660
$ws.Range("L116").Value = 3666.6667
$ws.Range("M116").Value = 782
$ws.Range("N116").Value = -10550.6667
$ws.Range("H135").Value = 71430440
$ws.Range("I135").Value = 31251714
$ws.Range("J135").Value = 200002370
$ws.Range("K135").Value = 281265426
$ws.Range("L135").Value = 1800021330
$ws.Range("M135").Value = -281262891
$ws.Range("N135").Value = -1800026400

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5335.8696
$ws.Range("I32").Value = 4166.7256
$ws.Range("J32").Value = 15691.143
$ws.Range("K32").Value = 4166.7256
$ws.Range("L32").Value = 15691.143
$ws.Range("M32").Value = -3879.7256
$ws.Range("N32").Value = -16265.143
$ws.Range("H61").Value = 4014.3513
$ws.Range("I61").Value = 3382.0386
$ws.Range("J61").Value = 5508.909
$ws.Range("K61").Value = 3382.0386
$ws.Range("L61").Value = 5508.909
$ws.Range("M61").Value = -3170.0386
$ws.Range("N61").Value = -5932.909
$ws.Range("H74").Value = 4735.273
$ws.Range("I74").Value = 1362.9032
$ws.Range("K74").Value = 1362.9032
$ws.Range("M74").Value = -488.9032
$ws.Range("H77").Value = 4735.273
$ws.Range("I77").Value = 1362.9032
$ws.Range("K77").Value = 6814.516
$ws.Range("M77").Value = -2446.516
$ws.Range("H132").Value = 6685.6665
$ws.Range("I132").Value = 2616.818
$ws.Range("J132").Value = 9483
$ws.Range("K132").Value = 7850.454000000001
$ws.Range("L132").Value = 28449
$ws.Range("M132").Value = -5320.454000000001
$ws.Range("N132").Value = -33509
$ws.Range("H136").Value = 4014.3513
$ws.Range("I136").Value = 3382.0386
$ws.Range("J136").Value = 5508.909
$ws.Range("K136").Value = 10146.1158
$ws.Range("L136").Value = 16526.727
$ws.Range("M136").Value = -7596.1158
$ws.Range("N136").Value = -21626.727

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 6360.115
$ws.Range("I134").Value = 7405.143
$ws.Range("K134").Value = 22215.429
$ws.Range("M134").Value = -19680.429

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4179.2856
$ws.Range("I31").Value = 3942.4375
$ws.Range("J31").Value = 5600.375
$ws.Range("K31").Value = 3942.4375
$ws.Range("L31").Value = 5600.375
$ws.Range("M31").Value = -3647.4375
$ws.Range("N31").Value = -6190.375
$ws.Range("H34").Value = 4179.2856
$ws.Range("I34").Value = 3942.4375
$ws.Range("J34").Value = 5600.375
$ws.Range("K34").Value = 3942.4375
$ws.Range("L34").Value = 5600.375
$ws.Range("M34").Value = -3740.4375
$ws.Range("N34").Value = -6004.375

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 73.666664
$ws.Range("I38").Value = 35
$ws.Range("J38").Value = 131.66667
$ws.Range("K38").Value = 105
$ws.Range("L38").Value = 395.00001
$ws.Range("M38").Value = 242
$ws.Range("N38").Value = -1089.00001

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 8933.333000000001
$ws.Range("I9").Value = 5500
$ws.Range("J9").Value = 15800
$ws.Range("K9").Value = 5500
$ws.Range("L9").Value = 15800
$ws.Range("M9").Value = -5330
$ws.Range("N9").Value = -16140
$ws.Range("H17").Value = 13223.111
$ws.Range("J17").Value = 19752
$ws.Range("L17").Value = 19752
$ws.Range("N17").Value = -20088
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("M19").ClearContents()
$ws.Range("H22").Value = 5750
$ws.Range("I22").Value = 5750
$ws.Range("K22").Value = 5750
$ws.Range("M22").Value = -5221
$ws.Range("H33").Value = 10000
$ws.Range("J33").Value = 10000
$ws.Range("L33").Value = 10000
$ws.Range("N33").Value = -10504
$ws.Range("H40").Value = 12018
$ws.Range("J40").Value = 12018
$ws.Range("L40").Value = 12018
$ws.Range("N40").Value = -12320
$ws.Range("H44").Value = 8000
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 8000
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 8000
$ws.Range("M44").ClearContents()
$ws.Range("N44").Value = -9192
$ws.Range("H46").Value = 29000
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 29000
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 29000
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -29312
$ws.Range("H57").Value = 20000
$ws.Range("I57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("M57").ClearContents()
$ws.Range("H58").Value = 49800
$ws.Range("I58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("M58").ClearContents()
$ws.Range("H102").Value = 6086.1875
$ws.Range("I102").Value = 5343.5454
$ws.Range("J102").Value = 7720
$ws.Range("K102").Value = 5343.5454
$ws.Range("L102").Value = 7720
$ws.Range("M102").Value = -3721.5454
$ws.Range("N102").Value = -10964
$ws.Range("H111").Value = 265000
$ws.Range("J111").Value = 265000
$ws.Range("L111").Value = 265000
$ws.Range("N111").Value = -271134
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()
$ws.Range("H132").Value = 1670.5518
$ws.Range("I132").Value = 1219.1177
$ws.Range("J132").Value = 2310.0833
$ws.Range("K132").Value = 3657.3531
$ws.Range("L132").Value = 6930.249899999999
$ws.Range("M132").Value = -1127.3531
$ws.Range("N132").Value = -11990.2499

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H36").Value = 68998
$ws.Range("J36").Value = 68998
$ws.Range("L36").Value = 68998
$ws.Range("N36").Value = -70122
$ws.Range("H55").Value = 182046.14
$ws.Range("I55").Value = 363784.62
$ws.Range("J55").Value = 307.63635
$ws.Range("K55").Value = 363784.62
$ws.Range("L55").Value = 307.63635
$ws.Range("M55").Value = -363611.62
$ws.Range("N55").Value = -653.63635
$ws.Range("H100").Value = 3081.5557
$ws.Range("I100").Value = 1750.1428
$ws.Range("J100").Value = 4515.385
$ws.Range("K100").Value = 1750.1428
$ws.Range("L100").Value = 4515.385
$ws.Range("M100").Value = -1209.1428
$ws.Range("N100").Value = -5597.385

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3521.3
$ws.Range("I62").Value = 3675.25
$ws.Range("J62").Value = 3418.6667
$ws.Range("K62").Value = 3675.25
$ws.Range("L62").Value = 3418.6667
$ws.Range("M62").Value = -3051.25
$ws.Range("N62").Value = -4666.6667
$ws.Range("H65").Value = 3521.3
$ws.Range("I65").Value = 3675.25
$ws.Range("J65").Value = 3418.6667
$ws.Range("K65").Value = 18376.25
$ws.Range("L65").Value = 17093.3335
$ws.Range("M65").Value = -15256.25
$ws.Range("N65").Value = -23333.3335
